$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.495.83"
$ws.Range("E2").Value = "  -6.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.246.37"
$ws.Range("E3").Value = "  -9.69%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "176.45"
$ws.Range("E5").Value = "  -12.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "508.37"
$ws.Range("E6").Value = "  -11.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  -5.47%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.247.16"
$ws.Range("E9").Value = "  -9.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.603"
$ws.Range("E10").Value = "  -11.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.66"
$ws.Range("E11").Value = "  -6.39%  "

$ws.Range("E12").Value = "  -13.35%  "

$ws.Range("E13").Value = "  -11.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.86"
$ws.Range("E14").Value = "  -14.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.767.87"
$ws.Range("E15").Value = "  -9.52%  "

$ws.Range("E16").Value = "  -4.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.247.64"
$ws.Range("E17").Value = "  -9.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.158.97"
$ws.Range("E18").Value = "  -6.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.82"
$ws.Range("E19").Value = "  -11.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.63"
$ws.Range("E20").Value = "  -13.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.926"
$ws.Range("E21").Value = "  -12.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "363.58"
$ws.Range("E22").Value = "  -9.87%  "

$ws.Range("E23").Value = "  -7.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.57"
$ws.Range("E24").Value = "  -15.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.60"

$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.71"
$ws.Range("E27").Value = "  -4.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.57"
$ws.Range("E28").Value = "  -11.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.92"
$ws.Range("E29").Value = "  -12.56%  "

$ws.Range("E30").Value = "  -11.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "635.17"
$ws.Range("E31").Value = "  -6.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.85"
$ws.Range("E32").Value = "  -11.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.46"
$ws.Range("E33").Value = "  -16.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.82"
$ws.Range("E34").Value = "  -10.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.26"
$ws.Range("E35").Value = "  -7.99%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.101"
$ws.Range("E37").Value = "  -11.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.05"
$ws.Range("E38").Value = "  -15.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.363"
$ws.Range("E39").Value = "  -11.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("E41").Value = "  -10.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.769.84"
$ws.Range("E42").Value = "  -13.17%  "

$ws.Range("E43").Value = "  -18.88%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0613"
$ws.Range("E44").Value = "  -19.92%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  -8.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0373"
$ws.Range("E46").Value = "  -9.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("E47").Value = "  -18.30%  "

$ws.Range("E48").Value = "  -7.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.08"
$ws.Range("E49").Value = "  -4.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  -5.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.74"
$ws.Range("E51").Value = "  -10.85%  "
